$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159 (shifts existing rows 159-251 down to 160-252)
$ws.Rows("159:159").Insert()

# Populate the newly inserted row 159 with the new weekly price record
$ws.Cells.Item(159, 1).Value = 4
$ws.Cells.Item(159, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(159, 3).Value = 'Los Lagos'
$ws.Cells.Item(159, 4).Value = 44606
$ws.Cells.Item(159, 5).Value = 10
$ws.Cells.Item(159, 6).Value = 100114014
$ws.Cells.Item(159, 7).Value = 'Betarraga'
$ws.Cells.Item(159, 8).Value = 'Sin especificar'
$ws.Cells.Item(159, 9).Value = 'Primera'
$ws.Cells.Item(159, 10).Value = 500
$ws.Cells.Item(159, 11).Value = 900
$ws.Cells.Item(159, 12).Value = 1000
$ws.Cells.Item(159, 13).Value = 950
$ws.Cells.Item(159, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(159, 15).Value = 'Región del Maule'
$ws.Cells.Item(159, 16).Value = 190
$ws.Cells.Item(159, 17).Value = 5
$ws.Cells.Item(159, 18).Value = 'Hortaliza'
